# Add the new match row (row 51) to the LaLiga 2023-2024 sheet, matching
# the existing table's layout/formatting (row 50 is the last data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold/border/alignment on Indice column, date number
# format on the match-date column, etc.) from the previous row so the new
# row renders identically to the rest of the table.
$ws.Range("A50:V50").Copy()
$ws.Range("A51:V51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row's values.
$ws.Range("A51").Value2 = 50
$ws.Range("B51").Value2 = "spain"
$ws.Range("C51").Value2 = "laliga"
$ws.Range("D51").Value2 = "2023-2024"
$ws.Range("E51").Value2 = 45191.875
$ws.Range("F51").Value2 = "Alaves"
$ws.Range("G51").Value2 = 0
$ws.Range("H51").Value2 = "Ath Bilbao"
$ws.Range("I51").Value2 = 2
$ws.Range("J51").Value2 = 3.83
$ws.Range("K51").Value2 = "11/09/2023 13:17"
$ws.Range("L51").Value2 = 4.25
$ws.Range("M51").Value2 = "22/09/2023 20:58"
$ws.Range("N51").Value2 = 3.11
$ws.Range("O51").Value2 = "11/09/2023 13:17"
$ws.Range("P51").Value2 = 3.1
$ws.Range("Q51").Value2 = "22/09/2023 20:58"
$ws.Range("R51").Value2 = 2.11
$ws.Range("S51").Value2 = "11/09/2023 13:17"
$ws.Range("T51").Value2 = 2.16
$ws.Range("U51").Value2 = "22/09/2023 20:58"
$ws.Range("V51").Value2 = "https://www.betexplorer.com/football/spain/laliga/alaves-ath-bilbao/p4l2FzZA/"
